# Apply the "new version with timestamp" edit:
#  - insert a new product row "FORFLOZIN 10MG 30 F.C. TABS" (before FUSI-ZON CREAM 15 GM)
#  - insert a new product row "TAVONIZA 20 MG 20 F.C.TABS." (before RICHI's successor, برفان القصاص)
#  - update the grand-total cell accordingly

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: make room for the two new rows -------------------------------
# Insert before the (old) row 11 - "برفان القصاص" - first (bottom-most insert first)
# so the row index for the first insertion below is unaffected by the later one.
$ws.Rows("11:11").Insert()
# Insert before the (old) row 6 - "FUSI-ZON CREAM 15 GM"
$ws.Rows("6:6").Insert()

# --- Step 2: clone formatting (fonts/fills/borders/number formats) onto ---
# --- the two freshly inserted (blank) rows, copying from a sibling row  ---
$ws.Range("A5:N5").Copy()
$ws.Range("A6:N6").PasteSpecial(-4122)

$ws.Range("A13:N13").Copy()
$ws.Range("A12:N12").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Step 3: merged cells for the two new rows -----------------------------
$ws.Range("B6:G6").Merge()
$ws.Range("H6:K6").Merge()
$ws.Range("L6:M6").Merge()

$ws.Range("B12:G12").Merge()
$ws.Range("H12:K12").Merge()
$ws.Range("L12:M12").Merge()

# --- Step 4: row heights (rows keep the row-height of whatever product ----
# --- ends up sitting on them, matching the canonical export)            ---
$ws.Rows(4).RowHeight  = 24.75
$ws.Rows(5).RowHeight  = 25.5
$ws.Rows(6).RowHeight  = 24.75
$ws.Rows(7).RowHeight  = 25.5
$ws.Rows(8).RowHeight  = 25.5
$ws.Rows(9).RowHeight  = 24.75
$ws.Rows(10).RowHeight = 25.5
$ws.Rows(11).RowHeight = 24.75
$ws.Rows(12).RowHeight = 25.5
$ws.Rows(13).RowHeight = 25.5
$ws.Rows(14).RowHeight = 24.75
$ws.Rows(15).RowHeight = 25.5
$ws.Rows(16).RowHeight = 24.75

# --- Step 5: cell values for the new FORFLOZIN row (row 6) ----------------
$ws.Range("A6").Value = 3
$ws.Range("B6").Value = "FORFLOZIN 10MG 30 F.C. TABS"
$ws.Range("H6").Value = "0:0"
$ws.Range("L6").Value = 102
$ws.Range("N6").Value = "0:0"

# --- Step 6: cell values for the new TAVONIZA row (row 12) ----------------
$ws.Range("A12").Value = 9
$ws.Range("B12").Value = "TAVONIZA 20 MG 20 F.C.TABS."
$ws.Range("H12").Value = "0:1"
$ws.Range("L12").Value = 99
$ws.Range("N12").Value = "0:2"

# --- Step 7: renumber the "م" (index) column for every row below --------
$ws.Range("A4").Value = 1
$ws.Range("A5").Value = 2
$ws.Range("A7").Value = 4
$ws.Range("A8").Value = 5
$ws.Range("A9").Value = 6
$ws.Range("A10").Value = 7
$ws.Range("A11").Value = 8
$ws.Range("A13").Value = 10
$ws.Range("A14").Value = 11
$ws.Range("A15").Value = 12
$ws.Range("A16").Value = 13

# --- Step 8: grand total (sum of the "سعر البيع" / L column) now includes -
# --- the prices of the two new rows ----------------------------------------
$ws.Range("K17").Value = 610.11

Write-Host "edit applied"
